# Regenerate save_data to use K (column G) instead of Strike# values.
# Update the "K" column (G2:G17) with newly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 6
    4  = 6
    5  = 5
    6  = 2
    7  = 1
    8  = 6
    9  = 1
    10 = 5
    11 = 3
    12 = 6
    13 = 4
    14 = 3
    15 = 1
    16 = 7
    17 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
